$d = $word.ActiveDocument
$d.Content.Find.Execute("target_altitude " + [char]8211 + " plane_altitude = 50cm.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "target_altitude " + [char]8211 + " plane_altitude = 17cm.", 2)
